$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Predictions")

# Rows 80-107: matches that have since completed. Update Status/Result/Resultado_Real/Profit/ROI/Enviado
$updates = @(
@(80, "Completed", "Away Win", "Acierto", 1.26, 105, "2025-09-14 04:25:49"),
    @(81, "Completed", "Home Win", "Acierto", 1.23, 95, "2025-09-14 04:25:49"),
    @(82, "Completed", "Home Win", "Acierto", 1.03, 115, "2025-09-14 04:25:49"),
    @(83, "Completed", "Away Win", "Acierto", 1.59, 55, "2025-09-14 04:25:49"),
    @(84, "Completed", "Draw", "Fallo", -1.5, -100, "2025-09-14 04:25:49"),
    @(85, "Completed", "Draw", "Fallo", -2.4, -100, "2025-09-14 04:25:49"),
    @(86, "Completed", "Home Win", "Acierto", 1.36, 91, "2025-09-14 04:25:49"),
    @(87, "Completed", "Home Win", "Acierto", 1.22, 42, "2025-09-14 04:25:49"),
    @(88, "Completed", "Away Win", "Fallo", -1.9, -100, "2025-09-14 04:25:49"),
    @(89, "Completed", "Away Win", "Acierto", 1.74, 60, "2025-09-14 04:25:49"),
    @(90, "Completed", "Home Win", "Acierto", 1.53, 85, "2025-09-14 04:25:49"),
    @(91, "Completed", "Home Win", "Acierto", 1.54, 70, "2025-09-14 04:25:49"),
    @(92, "Completed", "Away Win", "Acierto", 1.39, 48, "2025-09-14 04:25:49"),
    @(93, "Completed", "Away Win", "Acierto", 1.1, 110, "2025-09-14 04:25:49"),
    @(94, "Completed", "Away Win", "Fallo", -2.6, -100, "2025-09-14 04:25:49"),
    @(95, "Completed", "Home Win", "Acierto", 1.54, 53, "2025-09-14 04:25:49"),
    @(96, "Completed", "Away Win", "Acierto", 1.8, 62, "2025-09-14 04:25:49"),
    @(97, "Completed", "Home Win", "Acierto", 1.61, 70, "2025-09-14 04:25:49"),
    @(98, "Completed", "Away Win", "Fallo", -1.1, -100, "2025-09-14 04:25:49"),
    @(99, "Completed", "Away Win", "Acierto", 1.32, 120, "2025-09-14 04:25:49"),
    @(100, "Completed", "Home Win", "Fallo", -2.9, -100, "2025-09-14 15:19:33"),
    @(101, "Completed", "Draw", "Fallo", -2.3, -100, "2025-09-14 15:19:33"),
    @(102, "Completed", "Draw", "Fallo", -2.9, -100, "2025-09-14 15:19:33"),
    @(103, "Completed", "Home Win", "Acierto", 1.52, 80, "2025-09-14 15:19:33"),
    @(104, "Completed", "Away Win", "Acierto", 1.74, 60, "2025-09-14 15:19:33"),
    @(105, "Completed", "Away Win", "Acierto", 1.68, 70, "2025-09-14 15:19:33"),
    @(106, "Completed", "Draw", "Fallo", -0.9, -100, "2025-09-14 15:19:33"),
    @(107, "Completed", "Draw", "Fallo", -2.3, -100, "2025-09-14 15:19:33")
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 12).Value = $u[1]   # L: Status
    $ws.Cells.Item($row, 13).Value = $u[2]   # M: Result
    $ws.Cells.Item($row, 14).Value = $u[3]   # N: Resultado_Real
    $ws.Cells.Item($row, 15).Value = $u[4]   # O: Profit
    $ws.Cells.Item($row, 16).Value = $u[5]   # P: ROI
    $ws.Cells.Item($row, 17).Value = $u[6]   # Q: Enviado
}

# New rows 141-144: freshly synced predictions
$newRows = @(
@(141, "2025-09-15", "Liga de Expansión MX", "Irapuato", "Tepatitlán", "Home Win", "54.89%", 2.15, "16.83%", 1, 0.01566385563149848, 0.1566385563149848, "Pending"),
    @(142, "2025-09-15", "Süper Lig", "Rizespor", "Genclerbirligi", "Home Win", "77.51%", 1.62, "24.32%", 2.6, 0.04124871888302439, 0.4124871888302438, "Pending"),
    @(143, "2025-09-15", "Serie A", "Como", "Genoa", "Home Win", "77.08%", 1.65, "25.91%", 2.6, 0.04182086198902012, 0.4182086198902012, "Pending"),
    @(144, "2025-09-15", "La Liga", "Espanyol", "Mallorca", "Home Win", "53.12%", 2.2, "15.69%", 0.9, 0.01404630659222537, 0.1404630659222537, "Pending")
)

foreach ($n in $newRows) {
    $row = $n[0]
    $ws.Cells.Item($row, 1).Value = $n[1]    # A: Date
    $ws.Cells.Item($row, 2).Value = $n[2]    # B: Liga
    $ws.Cells.Item($row, 3).Value = $n[3]    # C: Local
    $ws.Cells.Item($row, 4).Value = $n[4]    # D: Visitante
    $ws.Cells.Item($row, 5).Value = $n[5]    # E: Prediccion
    $ws.Cells.Item($row, 6).Value = $n[6]    # F: Probabilidad
    $ws.Cells.Item($row, 7).Value = $n[7]    # G: Cuota_Bet365
    $ws.Cells.Item($row, 8).Value = $n[8]    # H: EV
    $ws.Cells.Item($row, 9).Value = $n[9]    # I: Stake
    $ws.Cells.Item($row, 10).Value = $n[10]  # J: StakePct
    $ws.Cells.Item($row, 11).Value = $n[11]  # K: KellyFrac
    $ws.Cells.Item($row, 12).Value = $n[12]  # L: Status
}
